# The "Total Wild Abundance" sheet had a stray title row ("TOTAL WILD
# ABUNDANCE") sitting above the real header row. Remove that row so the
# data starts with the column headers in row 1 (this also removes the
# now-unused "TOTAL WILD ABUNDANCE" shared string and shifts every
# subsequent row up by one, shrinking the used range from A1:V32 to
# A1:V31).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Total Wild Abundance")
$ws.Activate()

# Select the entire first row (mirrors the manual "click row header, then
# delete" workflow) before removing it.
$ws.Rows("1:1").Select()
$ws.Rows("1:1").Delete()
